$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 525
$ws.Range("J80").Value = 518.75
$ws.Range("L80").Value = 1556.25
$ws.Range("N80").Value = -3552.25
$ws.Range("H83").Value = 525
$ws.Range("J83").Value = 518.75
$ws.Range("L83").Value = 4668.75
$ws.Range("N83").Value = -14652.75
$ws.Range("H115").Value = 435
$ws.Range("I115").Value = 435
$ws.Range("K115").Value = 1305
$ws.Range("M115").Value = 262
$ws.Range("H116").Value = 2872.25
$ws.Range("I116").Value = 2663
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2663
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 779
$ws.Range("N116").Value = -10384
$ws.Range("H137").Value = 5854.8096
$ws.Range("I137").Value = 2665
$ws.Range("K137").Value = 7995
$ws.Range("M137").Value = -5445
$ws.Range("H141").Value = 5049.5
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 25000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 25000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 25000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -25860
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 37500
$ws.Range("I26").Value = 37500
$ws.Range("K26").Value = 37500
$ws.Range("M26").Value = -37208
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6626.125
$ws.Range("I94").Value = 3796.4
$ws.Range("K94").Value = 3796.4
$ws.Range("M94").Value = -3345.4
$ws.Range("H99").Value = 2343.318
$ws.Range("I99").Value = 1984.3334
$ws.Range("K99").Value = 1984.3334
$ws.Range("M99").Value = -486.3334
$ws.Range("H126").Value = 2343.318
$ws.Range("I126").Value = 1984.3334
$ws.Range("K126").Value = 5953.0002
$ws.Range("M126").Value = -3483.0002
$ws.Range("H132").Value = 729
$ws.Range("I132").Value = 673.75
$ws.Range("K132").Value = 2021.25
$ws.Range("M132").Value = 508.75
$ws.Range("H134").Value = 1909
$ws.Range("I134").Value = 1909
$ws.Range("K134").Value = 5727
$ws.Range("M134").Value = -3192
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 250
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 750
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1096
$ws.Range("H17").Value = 3794.45
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 4730.5625
$ws.Range("K17").Value = 150
$ws.Range("L17").Value = 14191.6875
$ws.Range("M17").Value = 19
$ws.Range("N17").Value = -14529.6875
$ws.Range("H68").Value = 1842.2858
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 1842.2858
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H80").Value = 4188.9287
$ws.Range("I80").Value = 3973.6316
$ws.Range("K80").Value = 11920.8948
$ws.Range("M80").Value = -10984.8948
$ws.Range("H83").Value = 4188.9287
$ws.Range("I83").Value = 3973.6316
$ws.Range("K83").Value = 35762.6844
$ws.Range("M83").Value = -31082.6844
$ws.Range("H86").Value = 379.2
$ws.Range("I86").Value = 386.875
$ws.Range("K86").Value = 1160.625
$ws.Range("M86").Value = 25.375
$ws.Range("H89").Value = 379.2
$ws.Range("I89").Value = 386.875
$ws.Range("K89").Value = 3481.875
$ws.Range("M89").Value = 2446.125
$ws.Range("H92").Value = 900
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 871.4666999999999
$ws.Range("I113").Value = 574.3333
$ws.Range("J113").Value = 945.75
$ws.Range("K113").Value = 1722.9999
$ws.Range("L113").Value = 2837.25
$ws.Range("M113").Value = 447.0001
$ws.Range("N113").Value = -7177.25
$ws.Range("H118").Value = 105
$ws.Range("I118").Value = 105
$ws.Range("K118").Value = 315
$ws.Range("M118").Value = 928
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 1612.9474
$ws.Range("I131").Value = 725.875
$ws.Range("K131").Value = 2177.625
$ws.Range("M131").Value = 2862.375
$ws.Range("H132").Value = 4699.2
$ws.Range("J132").Value = 5124.25
$ws.Range("L132").Value = 46118.25
$ws.Range("N132").Value = -51178.25
$ws.Range("H137").Value = 2870
$ws.Range("I137").Value = 1910
$ws.Range("J137").Value = 5750
$ws.Range("K137").Value = 5730
$ws.Range("L137").Value = 17250
$ws.Range("M137").Value = -630
$ws.Range("N137").Value = -27450
$ws.Range("H141").Value = 2207
$ws.Range("I141").Value = 2207
$ws.Range("K141").Value = 6621
$ws.Range("M141").Value = -1441
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 59831.65
$ws.Range("J20").Value = 500499
$ws.Range("L20").Value = 500499
$ws.Range("N20").Value = -500951
$ws.Range("H22").Value = 730.1667
$ws.Range("I22").Value = 718.4
$ws.Range("J22").Value = 789
$ws.Range("K22").Value = 718.4
$ws.Range("L22").Value = 789
$ws.Range("M22").Value = -423.4
$ws.Range("N22").Value = -1379
$ws.Range("H27").Value = 730.1667
$ws.Range("I27").Value = 718.4
$ws.Range("J27").Value = 789
$ws.Range("K27").Value = 718.4
$ws.Range("L27").Value = 789
$ws.Range("M27").Value = -611.4
$ws.Range("N27").Value = -1003
$ws.Range("H40").Value = 4599
$ws.Range("I40").Value = 4599
$ws.Range("K40").Value = 4599
$ws.Range("M40").Value = -4463
$ws.Range("H93").Value = 893.5
$ws.Range("I93").Value = 893.5
$ws.Range("K93").Value = 893.5
$ws.Range("M93").Value = 354.5
$ws.Range("H136").Value = 3741
$ws.Range("I136").Value = 2494.5
$ws.Range("K136").Value = 7483.5
$ws.Range("M136").Value = -4933.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6134.8237
$ws.Range("I126").Value = 4037.75
$ws.Range("K126").Value = 12113.25
$ws.Range("M126").Value = -9643.25
$ws.Range("H136").Value = 3233.6667
$ws.Range("I136").Value = 2063.125
$ws.Range("K136").Value = 6189.375
$ws.Range("M136").Value = -3639.375
